$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    Heading1 title paragraph ("Play Dragon Tribe Slot Game for
#    Free - Exciting Dragon Spins Feature") and before the
#    "DRAGON TRIBE SLOT GAME OVERVIEW" Heading2 paragraph.
#
#    We get a clean Normal-style blank paragraph (just <w:r/>, no
#    stray pPr/rsid noise) by splitting right before the first
#    existing Normal-style paragraph, then relocate the Heading2
#    paragraph after it so the ordering matches the target.
# ------------------------------------------------------------------

$firstNormalPara = $d.Paragraphs(3)          # "Dragon Tribe? Sounds like a bunch..."
$firstNormalPara.Range.InsertParagraphBefore() | Out-Null

$headingPara = $d.Paragraphs(2)              # "DRAGON TRIBE SLOT GAME OVERVIEW"
$headingPara.Range.Cut() | Out-Null

$metaPara = $d.Paragraphs(2)                 # now the freshly created blank paragraph
$afterMeta = $d.Range($metaPara.Range.End, $metaPara.Range.End)
$afterMeta.Paste() | Out-Null

# $metaPara is now positioned between the title and the (relocated)
# Heading2 paragraph. Fill it with the full text first (plain), then
# go back and bold only the "Meta description" label - applying the
# Bold formatting *after* the text exists is what keeps it isolated
# in its own run instead of bleeding into the text typed next to it.
$metaRange = $d.Paragraphs(2).Range
$insertPoint = $d.Range($metaRange.End - 1, $metaRange.End - 1)
$insertPoint.InsertBefore("Meta description: Experience the dragon-filled world of Dragon Tribe slots. Play for free and win up to 27,000x with the exciting Dragon Spins feature and xNudge Wilds.")

$metaRange = $d.Paragraphs(2).Range
$labelRange = $d.Range($metaRange.Start, $metaRange.Start + 16)
$labelRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Near the end of the document, remove the paragraph holding the
#    bold "Play Dragon Tribe Slot Game for Free - Exciting Dragon
#    Spins Feature" text entirely (it was duplicated to the top as
#    the meta-description paragraph above).
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Play Dragon Tribe Slot Game for Free - Exciting Dragon Spins Feature") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-prompt copy, keeping its italic formatting and leading
#    empty run untouched.
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$fullRange = $lastPara.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.Text = 'Create a feature image for Dragon Tribe, a cartoon-style image featuring a happy Maya warrior with glasses. The image should incorporate the dragon theme with flames and dragons flying in the background. Use bold and bright colors to make the image pop and catch the eye of potential players. Make sure to showcase the xNudge Wilds and Dragon Spins features in the image to give players a taste of the action-packed gameplay. The Maya warrior should be standing in front of the reel set, with the game''s logo at the top and the words "Dragon Tribe" written in a fun and playful font. Overall, the image should convey the excitement and adventure of this slot game.'

Write-Host "Done"
